# Apply the "Nature Green" template text-formatting fix:
#   - explicit left paragraph alignment (adds algn="l" on <a:pPr>)
#   - font typeface changed from Calibri to Arial
# on the two text boxes of slide 1 and the two text boxes of slide 2.

$p = $ppt.ActivePresentation

# ppAlignLeft
$ppAlignLeft = 1

function Set-LeftArial {
    param($shape)

    $tr = $shape.TextFrame.TextRange
    $tr.ParagraphFormat.Alignment = $ppAlignLeft
    $tr.Font.Name = "Arial"
}

# --- Slide 1 ---------------------------------------------------------
$s1 = $p.Slides.Item(1)

# TextBox 4 -> "Nature Green" (title, sz=4800 b=1, color 226422)
Set-LeftArial $s1.Shapes.Item(4)

# TextBox 5 -> "Eco-Friendly Template" (subtitle, sz=2000 b=0, color 507850)
Set-LeftArial $s1.Shapes.Item(5)

# --- Slide 2 ---------------------------------------------------------
$s2 = $p.Slides.Item(2)

# TextBox 3 -> "Eco Features" (title, sz=3200 b=1, color FFFFFF)
Set-LeftArial $s2.Shapes.Item(3)

# TextBox 5 -> bullet list (sz=2000 b=0, color 325032)
Set-LeftArial $s2.Shapes.Item(5)
